# Updates the "cryptos" price table (generated from the upstream GitHub
# Actions data refresh) to the new Price (column D) / Volume(1h) (column E)
# figures, matching coin names/links in rows 45-46 (Aptos and BabyDogeCoin
# swapped places in the source feed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D sometimes holds price strings that look numeric (e.g. "240.26").
# Assigning those directly would let Excel auto-convert the cell to a real
# number, which would not match the source sheet where every Price/Volume
# cell is stored as text. Prefixing the value with an apostrophe forces a
# text entry; resetting the style afterward drops the leftover quote-prefix
# formatting so the cell ends up with no explicit style - same as in the
# original file.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2: Bitcoin
$ws.Range('D2').Value = '29.425.48'
$ws.Range('E2').Value = '  +0.12%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.850.30'
$ws.Range('E3').Value = '  +0.10%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.04%  '

# Row 5: BNB
Set-TextValue $ws.Range('D5') '240.26'
$ws.Range('E5').Value = '  +0.04%  '

# Row 6: XRP
Set-TextValue $ws.Range('D6') '0.6291'
$ws.Range('E6').Value = '  -0.03%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.02%  '

# Row 8: Dogecoin
Set-TextValue $ws.Range('D8') '0.07649'
$ws.Range('E8').Value = '  +0.75%  '

# Row 9: Cardano
$ws.Range('E9').Value = '  -0.53%  '

# Row 10: Solana
Set-TextValue $ws.Range('D10') '24.90'
$ws.Range('E10').Value = '  +1.91%  '

# Row 11: WrappedEther
$ws.Range('D11').Value = '2.111.72'
$ws.Range('E11').Value = '  +13.99%  '

# Row 12: TRON
Set-TextValue $ws.Range('D12') '0.07738'
$ws.Range('E12').Value = '  -0.02%  '

# Row 13: Polkadot
Set-TextValue $ws.Range('D13') '5.033'

# Row 14: Polygon
Set-TextValue $ws.Range('D14') '0.6811'
$ws.Range('E14').Value = '  +0.45%  '

# Row 15: ShibaInu
Set-TextValue $ws.Range('D15') '0.00001064'
$ws.Range('E15').Value = '  -1.26%  '

# Row 17: Uniswap
Set-TextValue $ws.Range('D17') '6.184'
$ws.Range('E17').Value = '  +0.20%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '29.487.27'
$ws.Range('E18').Value = '  +0.30%  '

# Row 19: BitcoinCash
Set-TextValue $ws.Range('D19') '228.54'
$ws.Range('E19').Value = '  +0.09%  '

# Row 20: Avalanche
$ws.Range('E20').Value = '  -0.88%  '

# Row 21: Dai
$ws.Range('E21').Value = '  +0.03%  '

# Row 22: Chainlink
Set-TextValue $ws.Range('D22') '7.465'
$ws.Range('E22').Value = '  -0.18%  '

# Row 23: BinanceUSD
Set-TextValue $ws.Range('D23') '1.001'
$ws.Range('E23').Value = '  -0.04%  '

# Row 24: Monero
Set-TextValue $ws.Range('D24') '157.56'
$ws.Range('E24').Value = '  +0.20%  '

# Row 25: Stellar
Set-TextValue $ws.Range('D25') '0.1385'
$ws.Range('E25').Value = '  -0.68%  '

# Row 26: Cosmos
Set-TextValue $ws.Range('D26') '8.434'

# Row 27: EthereumClassic
Set-TextValue $ws.Range('D27') '17.69'
$ws.Range('E27').Value = '  +0.43%  '

# Row 28: Toncoin
Set-TextValue $ws.Range('D28') '1.387'
$ws.Range('E28').Value = '  +6.67%  '

# Row 29: PancakeSwap
$ws.Range('E29').Value = '  -0.24%  '

# Row 30: Hedera
Set-TextValue $ws.Range('D30') '0.05613'
$ws.Range('E30').Value = '  +0.48%  '

# Row 31: Filecoin
Set-TextValue $ws.Range('D31') '4.133'

# Row 32: InternetComputer(DFINITY)
Set-TextValue $ws.Range('D32') '4.052'
$ws.Range('E32').Value = '  +0.63%  '

# Row 33: LidoDAOToken
Set-TextValue $ws.Range('D33') '1.845'
$ws.Range('E33').Value = '  +0.21%  '

# Row 34: ARBITRUM
$ws.Range('E34').Value = '  +0.73%  '

# Row 35: ImmutableX
Set-TextValue $ws.Range('D35') '0.7016'
$ws.Range('E35').Value = '  -1.11%  '

# Row 36: HuobiToken
Set-TextValue $ws.Range('D36') '2.592'
$ws.Range('E36').Value = '  +0.19%  '

# Row 37: VeChain
Set-TextValue $ws.Range('D37') '0.01803'
$ws.Range('E37').Value = '  +0.20%  '

# Row 38: Maker
$ws.Range('D38').Value = '1.229.82'
$ws.Range('E38').Value = '  -0.66%  '

# Row 39: MXToken
Set-TextValue $ws.Range('D39') '2.711'
$ws.Range('E39').Value = '  -2.14%  '

# Row 40: FraxShare
Set-TextValue $ws.Range('D40') '6.442'
$ws.Range('E40').Value = '  +0.48%  '

# Row 41: TrustWalletToken
Set-TextValue $ws.Range('D41') '0.9071'
$ws.Range('E41').Value = '  +0.24%  '

# Row 42: PaxDollar
Set-TextValue $ws.Range('D42') '1.001'
$ws.Range('E42').Value = '  +0.06%  '

# Row 43: Quant
Set-TextValue $ws.Range('D43') '102.33'
$ws.Range('E43').Value = '  +0.69%  '

# Row 44: Aave
Set-TextValue $ws.Range('D44') '66.09'
$ws.Range('E44').Value = '  +0.21%  '

# Row 45: Aptos
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D45') '0.00000000121'
$ws.Range('E45').Value = '  -0.09%  '

# Row 46: BabyDogeCoin
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D46') '7.195'
$ws.Range('E46').Value = '  +0.71%  '

# Row 47: TheSandbox
$ws.Range('E47').Value = '  +0.22%  '

# Row 48: Algorand
Set-TextValue $ws.Range('D48') '0.1156'

# Row 49: EnergySwap
Set-TextValue $ws.Range('D49') '9.002'
$ws.Range('E49').Value = '  -0.21%  '

# Row 50: RenderToken
Set-TextValue $ws.Range('D50') '1.685'
$ws.Range('E50').Value = '  +0.58%  '

# Row 51: Cronos
$ws.Range('E51').Value = '  -0.06%  '
